$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "services": append two new service rows (17 -> id 16, 18 -> id 17)
# ---------------------------------------------------------------------------
$wsServices = $wb.Worksheets.Item("services")

$wsServices.Range("A17").Value = 16
$wsServices.Range("B17").Value = "faxina"
$wsServices.Range("D17").Value = 47.4
$wsServices.Range("G17").Value = "unidade"
$wsServices.Range("M17").Value = "2025-09-18T21:24:11.049955"
$wsServices.Range("N17").Value = "2025-09-18T21:24:11.049955"

$wsServices.Range("A18").Value = 17
$wsServices.Range("B18").Value = "eletricista"
$wsServices.Range("D18").Value = 47.4
$wsServices.Range("G18").Value = "unidade"
$wsServices.Range("M18").Value = "2025-09-18T21:27:57.779088"
$wsServices.Range("N18").Value = "2025-09-18T21:27:57.779088"

# ---------------------------------------------------------------------------
# Sheet "quotes": append two new quote rows (11 -> id 11, 12 -> id 12)
# ---------------------------------------------------------------------------
$wsQuotes = $wb.Worksheets.Item("quotes")

$wsQuotes.Range("A11").Value = 11
$wsQuotes.Range("B11").Value = "ORC202509009"
$wsQuotes.Range("C11").Value = 1
$wsQuotes.Range("D11").Value = "Orçamento - faxina"
$wsQuotes.Range("E11").Value = "limpeza"
$wsQuotes.Range("H11").Value = "pendente"
$wsQuotes.Range("M11").Value = 47.4
$wsQuotes.Range("R11").Value = "2025-09-18T21:24:11.937131"
$wsQuotes.Range("S11").Value = "2025-09-18T21:24:11.937131"

$wsQuotes.Range("A12").Value = 12
$wsQuotes.Range("B12").Value = "ORC202509010"
$wsQuotes.Range("C12").Value = 1
$wsQuotes.Range("D12").Value = "Orçamento - eletricista"
$wsQuotes.Range("E12").Value = "manutencao"
$wsQuotes.Range("H12").Value = "pendente"
$wsQuotes.Range("M12").Value = 47.4
$wsQuotes.Range("R12").Value = "2025-09-18T21:27:58.411126"
$wsQuotes.Range("S12").Value = "2025-09-18T21:27:58.411126"

# ---------------------------------------------------------------------------
# Sheet "quote_items": append two new quote-item rows (12 -> id 13, 13 -> id 14)
# ---------------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("quote_items")

$wsItems.Range("A12").Value = 13
$wsItems.Range("B12").Value = 11
$wsItems.Range("C12").Value = 16
$wsItems.Range("D12").Value = 1
$wsItems.Range("E12").Value = 47.4
$wsItems.Range("G12").Value = 47.4
$wsItems.Range("H12").Value = "faxina"
$wsItems.Range("J12").Value = "unidade"
$wsItems.Range("O12").Value = "2025-09-18T21:24:11.937131"

$wsItems.Range("A13").Value = 14
$wsItems.Range("B13").Value = 12
$wsItems.Range("C13").Value = 17
$wsItems.Range("D13").Value = 1
$wsItems.Range("E13").Value = 47.4
$wsItems.Range("G13").Value = 47.4
$wsItems.Range("H13").Value = "eletricista"
$wsItems.Range("J13").Value = "unidade"
$wsItems.Range("O13").Value = "2025-09-18T21:27:58.411126"
